$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New diary entry for Jan 23rd, 2020 goes into the currently-blank
# --- template row (row 12); the blank template itself is pushed down
# --- into row 14 (reusing the placeholder text that used to live in row 12).

# Row 12: fill in the new diary entry.
$ws.Range("A12").Value = "Jan 23rd, 2020"
$ws.Range("B12").Value = "5 - 8 pm"
$ws.Range("C12").Value = "N/A"
$ws.Range("D12").Value = "Revise last week's material, learn mental models and externalizing them, get familiar with UML class diagrams, do an in-class practice and have a face-to-face communication with Alegria."
$ws.Range("E12").Value = "Knew what mental models and UML class diagrams are and gave a try to make UML diagrams with the simpleUML plugin. What's more, learned about what a programmer's life (basically her way to read code) is like from Alegria's narrate."
$ws.Range("F12").Value = "UML class diagrams help a lot when programmers try to figure out logical relationships among many java components by visualize them clearly. "
$ws.Range("G12").Value = "It's about the time to devote more time to our group project and really dive into it."

# Row 12 grows taller to fit the new wrapped text.
$ws.Rows.Item(12).RowHeight = 78

# Row 14: restore the blank-template placeholder prompts (previously in row 12).
$ws.Range("A14").Value = "<what day?>"
$ws.Range("B14").Value = "<what time?>"
$ws.Range("C14").Value = "<as applicable, with whom?>"
$ws.Range("D14").Value = "<what did you want to accomplish?>"
$ws.Range("E14").Value = "<what did you actually accomplish?>"
$ws.Range("F14").Value = "<what insight(s) did you gain?>"
$ws.Range("G14").Value = "<how did you feel during the activity?>"

# Move the active selection/view to reflect where the author was working.
$null = $ws.Range("F12").Select()
